$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a new worksheet "2022-Q4" right before the current "2022-Q3"
#    sheet (index 2). All later quarter sheets shift right by one position
#    automatically and keep their original content untouched.
# ---------------------------------------------------------------------------
$anchor = $wb.Worksheets.Item(2)            # "2022-Q3" (before insertion)
$q4 = $wb.Worksheets.Add($anchor)
$q4.Name = "2022-Q4"

# "2022-Q3" is now one slot to the right of where it used to be; grab it so
# we can clone its layout/formatting (headers + column A style) onto the
# brand-new sheet.
$q3 = $wb.Worksheets.Item(3)
$q3.Range("A1:H3").Copy()
$q4.Range("A1").PasteSpecial(-4122)          # xlPasteFormats
$q4.Range("A1").ClearContents()

# ---------------------------------------------------------------------------
# 2. Fill the new "2022-Q4" sheet with its data (same two funds as the
#    "2022-Q3" report, refreshed numbers for the new quarter).
# ---------------------------------------------------------------------------
$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"

$q4.Range("A2").Value = 0
$q4.Range("B2").NumberFormat = "@"
$q4.Range("B2").Value = "519956"
$q4.Range("C2").Value = "长信睿进灵活配置混合C"
$q4.Range("D2").NumberFormat = "@"
$q4.Range("D2").Value = "6.44"
$q4.Range("E2").NumberFormat = "@"
$q4.Range("E2").Value = "43.59"
$q4.Range("F2").NumberFormat = "@"
$q4.Range("F2").Value = "2.80"
$q4.Range("G2").NumberFormat = "@"
$q4.Range("G2").Value = "0.1803"
$q4.Range("H2").Value = 5

$q4.Range("A3").Value = 1
$q4.Range("B3").NumberFormat = "@"
$q4.Range("B3").Value = "519957"
$q4.Range("C3").Value = "长信睿进灵活配置混合A"
$q4.Range("D3").NumberFormat = "@"
$q4.Range("D3").Value = "0.01"
$q4.Range("E3").NumberFormat = "@"
$q4.Range("E3").Value = "43.59"
$q4.Range("F3").NumberFormat = "@"
$q4.Range("F3").Value = "2.80"
$q4.Range("G3").NumberFormat = "@"
$q4.Range("G3").Value = "0.0003"
$q4.Range("H3").Value = 5

# ---------------------------------------------------------------------------
# 3. Update the "总计" (summary) sheet: push the existing quarter rows down
#    by one and insert a new row for 2022-Q4 at the top.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)

$oldB2 = $total.Range("B2").Value()
$oldC2 = $total.Range("C2").Value()
$oldD2 = $total.Range("D2").Value()
$oldB3 = $total.Range("B3").Value()
$oldC3 = $total.Range("C3").Value()
$oldD3 = $total.Range("D3").Value()
$oldB4 = $total.Range("B4").Value()
$oldC4 = $total.Range("C4").Value()
$oldD4 = $total.Range("D4").Value()
$oldB5 = $total.Range("B5").Value()
$oldC5 = $total.Range("C5").Value()
$oldD5 = $total.Range("D5").Value()

# shift rows 2-5 down to rows 3-6 (bottom-up so we never clobber unread data)
$total.Range("B6").Value = $oldB5
$total.Range("C6").Value = $oldC5
$total.Range("D6").Value = $oldD5

$total.Range("B5").Value = $oldB4
$total.Range("C5").Value = $oldC4
$total.Range("D5").Value = $oldD4

$total.Range("B4").Value = $oldB3
$total.Range("C4").Value = $oldC3
$total.Range("D4").Value = $oldD3

$total.Range("B3").Value = $oldB2
$total.Range("C3").Value = $oldC2
$total.Range("D3").Value = $oldD2

# new row 6 needs the numbered-row style (column A) copied down
$total.Range("A5").Copy()
$total.Range("A6").PasteSpecial(-4122)
$total.Range("A6").Value = 4

# row 2 becomes the new 2022-Q4 entry
$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0.18

Write-Host "edit complete"
